# Crotone_matches_2020.xlsx update
# - Rows 2-9: refresh the xG_away / goals_home / goals_away values (columns E, F, G)
# - Rows 10-15: populate the previously-empty xG_home / xG_away / goals_home / goals_away
#   values (columns D, E, F, G)
#
# The source data stores every value (including the numeric-looking ones such as
# "4", "0", "1.80058", ...) as TEXT (shared strings), matching how the rest of the
# sheet was authored. A plain `Range.Value = "4"` assignment would make Excel treat
# the text as a number, so instead we render the literal text via TEXT(...,"General")
# in a scratch cell and paste-special the *value* into the destination cell. That
# keeps the destination's number format/style untouched (no stray styles) while still
# landing a genuine text cell.
#
# Cells are written column-by-column (D, then E, then F, then G) to mirror the
# original column-major export order of this dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell safely outside the sheet's real data (A1:G39).
$scratch = $ws.Range("ZZ1")

function Set-TextValue($cellRef, $val) {
    $scratch.Formula = "=TEXT(" + $val + ",""General"")"
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

# ---- Column D: rows 10-15 (new cells; D2:D9 stay as-is) ----
Set-TextValue "D10" "1.80058"
Set-TextValue "D11" "0.508882"
Set-TextValue "D12" "1.69279"
Set-TextValue "D13" "1.10391"
Set-TextValue "D14" "2.48022"
Set-TextValue "D15" "0.7999"

# ---- Column E: rows 2-9 updated, rows 10-15 newly added ----
Set-TextValue "E2" "1.32278"
Set-TextValue "E3" "2.43073"
Set-TextValue "E4" "1.5945"
Set-TextValue "E5" "1.58494"
Set-TextValue "E6" "1.40714"
Set-TextValue "E7" "3.49595"
Set-TextValue "E8" "0.300792"
Set-TextValue "E9" "1.79345"
Set-TextValue "E10" "0.401219"
Set-TextValue "E11" "1.93801"
Set-TextValue "E12" "0.824038"
Set-TextValue "E13" "0.123143"
Set-TextValue "E14" "1.04297"
Set-TextValue "E15" "1.50226"

# ---- Column F: rows 2-9 updated, rows 10-15 newly added ----
Set-TextValue "F2" "4"
Set-TextValue "F3" "0"
Set-TextValue "F4" "4"
Set-TextValue "F5" "1"
Set-TextValue "F6" "4"
Set-TextValue "F7" "1"
Set-TextValue "F8" "0"
Set-TextValue "F9" "0"
Set-TextValue "F10" "1"
Set-TextValue "F11" "0"
Set-TextValue "F12" "4"
Set-TextValue "F13" "0"
Set-TextValue "F14" "3"
Set-TextValue "F15" "2"

# ---- Column G: rows 2-9 updated, rows 10-15 newly added ----
Set-TextValue "G2" "1"
Set-TextValue "G3" "2"
Set-TextValue "G4" "1"
Set-TextValue "G5" "1"
Set-TextValue "G6" "2"
Set-TextValue "G7" "2"
Set-TextValue "G8" "0"
Set-TextValue "G9" "2"
Set-TextValue "G10" "0"
Set-TextValue "G11" "4"
Set-TextValue "G12" "1"
Set-TextValue "G13" "0"
Set-TextValue "G14" "1"
Set-TextValue "G15" "1"

# Clean up the scratch cell and clipboard marquee.
$scratch.ClearContents()
$excel.CutCopyMode = 0
